$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '68.251.51'
    'E2' = '  -1.88%  '
    'D3' = '2.450.83'
    'E3' = '  -1.87%  '
    'E4' = '  -0.02%  '
    'D5' = '558.51'
    'E5' = '  -2.90%  '
    'D6' = '163.00'
    'E6' = '  -2.22%  '
    'E7' = '  -0.01%  '
    'E8' = '  -1.95%  '
    'D9' = '2.449.67'
    'E9' = '  -1.84%  '
    'D10' = '0.150'
    'E10' = '  -5.66%  '
    'E11' = '  -1.95%  '
    'E12' = '  -5.36%  '
    'E13' = '  -2.62%  '
    'D14' = '2.909.06'
    'E14' = '  -1.63%  '
    'D15' = '68.288.63'
    'E15' = '  -1.64%  '
    'E16' = '  -3.53%  '
    'D17' = '23.24'
    'E17' = '  -5.93%  '
    'D18' = '2.508.75'
    'E18' = '  +0.36%  '
    'D19' = '10.97'
    'E19' = '  -2.23%  '
    'D20' = '7.17'
    'E20' = '  -3.88%  '
    'D21' = '341.86'
    'E21' = '  -1.71%  '
    'E22' = '  -3.37%  '
    'E23' = '  -0.24%  '
    'E24' = '  -3.91%  '
    'D25' = '67.64'
    'E25' = '  -4.27%  '
    'E26' = '  +7.24%  '
    'D27' = '3.70'
    'E27' = '  -6.32%  '
    'E28' = '  -1.48%  '
    'D29' = '8.14'
    'E29' = '  -6.99%  '
    'E30' = '  -6.41%  '
    'E31' = '  -7.64%  '
    'E32' = '  +129.95%  '
    'E33' = '  +0.05%  '
    'D34' = '432.16'
    'E34' = '  -4.98%  '
    'E35' = '  -3.54%  '
    'E36' = '  -3.52%  '
    'D37' = '156.93'
    'E37' = '  -0.35%  '
    'D38' = '19.00'
    'E38' = '  -0.30%  '
    'E39' = '  -0.03%  '
    'E40' = '  -5.64%  '
    'E41' = '  -3.06%  '
    'E42' = '  -3.50%  '
    'E43' = '  -4.93%  '
    'E44' = '  -5.22%  '
    'E45' = '  +0.22%  '
    'E46' = '  -6.26%  '
    'D47' = '133.65'
    'E47' = '  -5.42%  '
    'E48' = '  -3.54%  '
    'E49' = '  -1.95%  '
    'E50' = '  -6.81%  '
    'E51' = '  -3.06%  '
}

foreach ($cell in $updates.Keys) {
    # Force text interpretation so values like "163.00" / "0.150" / "19.00"
    # keep their literal digits instead of being coerced into numbers (which
    # would drop the trailing zeros). Restore the default "Normal" style
    # afterwards so no stray style index is left on the cell.
    $ws.Range($cell).NumberFormat = '@'
    $ws.Range($cell).Value = $updates[$cell]
    $ws.Range($cell).Style = 'Normal'
}